# Guidance of Reproducing the paper.pptx -- apply commit "uploaded updates to codes and folder"
$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Move slide 17 ("Extract Strongest Context Indicator") up to position 16,
#    pushing the former slide 16 ("Semantic Similarity") down to position 17.
# ---------------------------------------------------------------------------
$moved = $p.Slides.Item(17)
$moved.MoveTo(16)

# ---------------------------------------------------------------------------
# 2) Update the title of the slide that now sits at position 16 so it reads
#    "...As the Definition".
# ---------------------------------------------------------------------------
$slide16 = $p.Slides.Item(16)
$title16 = $slide16.Shapes.Title
$title16.TextFrame.TextRange.Text = "Extract Strongest Context Indicator As the Definition"

# ---------------------------------------------------------------------------
# 3) Fix the typo "patter" -> "pattern" on slide 7.
# ---------------------------------------------------------------------------
$slide7 = $p.Slides.Item(7)
$body7 = $slide7.Shapes.Item(2)
$tr7 = $body7.TextFrame.TextRange
$para7 = $tr7.Paragraphs(5)
$fullPara7 = $tr7.Characters($para7.Start, $para7.Length)
$fullPara7.Text = "The proportion of transactions containing the pattern in the entire transaction dataset"

# ---------------------------------------------------------------------------
# 4) Slide 21 text updates: generalize "given author" -> "given pattern"
#    wording on the annotation diagram.
# ---------------------------------------------------------------------------
$slide21 = $p.Slides.Item(21)

$rect11 = $slide21.Shapes.Item(10)
$rect11.TextFrame.TextRange.Text = "Annotation of Given Pattern by Context Units with Top Weights"

$rect12 = $slide21.Shapes.Item(11)
$rect12.TextFrame.TextRange.Text = "Find representative titles of given pattern"

$rect13 = $slide21.Shapes.Item(12)
$tr13 = $rect13.TextFrame.TextRange
$splitAt = "Find synonyms of given author ".IndexOf("given author ") + 1
$part2 = $tr13.Characters($splitAt, "given author ".Length)
$part2.Text = "given pattern "

# ---------------------------------------------------------------------------
# 5) Refresh the cached "last printed" date field on the slide master and on
#    every slide layout (12/10/2020 -> 12/12/2020).
# ---------------------------------------------------------------------------
function Update-CachedDate {
    param($shapes)
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "12/10/2020") {
                $whole = $tr.Characters(1, $tr.Length)
                $whole.Text = "12/12/2020"
            }
        }
    }
}

$master = $p.SlideMaster
Update-CachedDate $master.Shapes
for ($j = 1; $j -le $master.CustomLayouts.Count; $j++) {
    $layout = $master.CustomLayouts.Item($j)
    Update-CachedDate $layout.Shapes
}
